# Add chat bot message log rows 14-27 (category selection, group choice, and 5 tasks sent)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$c = $ws.Cells.Item(14, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 2)
$c.NumberFormat = '@'
$c.Value = '14:00:31'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 5)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 6)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 8)
$c.NumberFormat = '@'
$c.Value = 'a,dls,dfs,d;l,x''sd,c'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 9)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 10)
$c.NumberFormat = '@'
$c.Value = 'inknkn]'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 11)
$c.NumberFormat = '@'
$c.Value = 'jnkn;n'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 12)
$c.NumberFormat = '@'
$c.Value = 'jknl;'
$c.ClearFormats()
$c = $ws.Cells.Item(14, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #1'
$c.ClearFormats()

# Row 15
$c = $ws.Cells.Item(15, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 2)
$c.NumberFormat = '@'
$c.Value = '14:05:23'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 5)
$c.NumberFormat = '@'
$c.Value = 'Институт'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 6)
$c.NumberFormat = '@'
$c.Value = 'ВУЗ'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 8)
$c.NumberFormat = '@'
$c.Value = 'XcˀṢfd'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 9)
$c.NumberFormat = '@'
$c.Value = 'dfvbsdv'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 10)
$c.NumberFormat = '@'
$c.Value = 'sfbvxvc'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 11)
$c.NumberFormat = '@'
$c.Value = 'sfvxcfv'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 12)
$c.NumberFormat = '@'
$c.Value = 'xdfvxdfv'
$c.ClearFormats()
$c = $ws.Cells.Item(15, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #1'
$c.ClearFormats()

# Row 16
$c = $ws.Cells.Item(16, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 2)
$c.NumberFormat = '@'
$c.Value = '14:22:15'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 6)
$c.NumberFormat = '@'
$c.Value = 'Group B'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 7)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 8)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 9)
$c.NumberFormat = '@'
$c.Value = 'w;d,fs,d'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 10)
$c.NumberFormat = '@'
$c.Value = 'sdfcvsdv'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 11)
$c.NumberFormat = '@'
$c.Value = 'sdvsdv'
$c.ClearFormats()
$c = $ws.Cells.Item(16, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #3'
$c.ClearFormats()

# Row 17
$c = $ws.Cells.Item(17, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 2)
$c.NumberFormat = '@'
$c.Value = '14:25:09'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 6)
$c.NumberFormat = '@'
$c.Value = 'Group B'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 8)
$c.NumberFormat = '@'
$c.Value = 'ijojojo'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 9)
$c.NumberFormat = '@'
$c.Value = 'ipjlk'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 10)
$c.NumberFormat = '@'
$c.Value = 'uobjb'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 11)
$c.NumberFormat = '@'
$c.Value = 'ibkb'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 12)
$c.NumberFormat = '@'
$c.Value = '13000'
$c.ClearFormats()
$c = $ws.Cells.Item(17, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #3'
$c.ClearFormats()

# Row 18
$c = $ws.Cells.Item(18, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 2)
$c.NumberFormat = '@'
$c.Value = '14:26:43'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 6)
$c.NumberFormat = '@'
$c.Value = 'Group B'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 8)
$c.NumberFormat = '@'
$c.Value = 'Дддд'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 9)
$c.NumberFormat = '@'
$c.Value = 'Ддд'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 10)
$c.NumberFormat = '@'
$c.Value = 'Ддд'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 11)
$c.NumberFormat = '@'
$c.Value = 'Ддд'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 12)
$c.NumberFormat = '@'
$c.Value = 'Ддд'
$c.ClearFormats()
$c = $ws.Cells.Item(18, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #2'
$c.ClearFormats()

# Row 19
$c = $ws.Cells.Item(19, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 2)
$c.NumberFormat = '@'
$c.Value = '14:28:42'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 6)
$c.NumberFormat = '@'
$c.Value = 'Group C'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 7)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 8)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 9)
$c.NumberFormat = '@'
$c.Value = 'Щщз'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 10)
$c.NumberFormat = '@'
$c.Value = 'Лл'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 11)
$c.NumberFormat = '@'
$c.Value = 'Лл'
$c.ClearFormats()
$c = $ws.Cells.Item(19, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #3'
$c.ClearFormats()

# Row 20
$c = $ws.Cells.Item(20, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 2)
$c.NumberFormat = '@'
$c.Value = '14:30:08'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 6)
$c.NumberFormat = '@'
$c.Value = 'Group C'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 8)
$c.NumberFormat = '@'
$c.Value = 'Зз'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 9)
$c.NumberFormat = '@'
$c.Value = 'Зз'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 10)
$c.NumberFormat = '@'
$c.Value = 'Зз'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 11)
$c.NumberFormat = '@'
$c.Value = 'Дд'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 12)
$c.NumberFormat = '@'
$c.Value = 'Дщз'
$c.ClearFormats()
$c = $ws.Cells.Item(20, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #2'
$c.ClearFormats()

# Row 21
$c = $ws.Cells.Item(21, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 2)
$c.NumberFormat = '@'
$c.Value = '14:31:50'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 6)
$c.NumberFormat = '@'
$c.Value = 'Group A'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 8)
$c.NumberFormat = '@'
$c.Value = 'Цдщцщк'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 9)
$c.NumberFormat = '@'
$c.Value = 'Дцдадп'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 10)
$c.NumberFormat = '@'
$c.Value = 'Цллапш'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 11)
$c.NumberFormat = '@'
$c.Value = 'Дцла'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 12)
$c.NumberFormat = '@'
$c.Value = 'Лцлал'
$c.ClearFormats()
$c = $ws.Cells.Item(21, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #3'
$c.ClearFormats()

# Row 22
$c = $ws.Cells.Item(22, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 2)
$c.NumberFormat = '@'
$c.Value = '14:35:06'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 6)
$c.NumberFormat = '@'
$c.Value = 'Group B'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 7)
$c.NumberFormat = '@'
$c.Value = 'Нет'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 8)
$c.NumberFormat = '@'
$c.Value = 'Щщщщ'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 9)
$c.NumberFormat = '@'
$c.Value = 'Дддд'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 10)
$c.NumberFormat = '@'
$c.Value = 'Эээ'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 11)
$c.NumberFormat = '@'
$c.Value = 'Ввв'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 12)
$c.NumberFormat = '@'
$c.Value = 'Миол'
$c.ClearFormats()
$c = $ws.Cells.Item(22, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #2'
$c.ClearFormats()

# Row 23
$c = $ws.Cells.Item(23, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 2)
$c.NumberFormat = '@'
$c.Value = '15:04:53'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 6)
$c.NumberFormat = '@'
$c.Value = 'Group C'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 7)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 8)
$c.NumberFormat = '@'
$c.Value = 'Щщлл'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 9)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 10)
$c.NumberFormat = '@'
$c.Value = 'Дщз'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 11)
$c.NumberFormat = '@'
$c.Value = 'Жз'
$c.ClearFormats()
$c = $ws.Cells.Item(23, 13)
$c.NumberFormat = '@'
$c.Value = 'Ьл'
$c.ClearFormats()

# Row 24
$c = $ws.Cells.Item(24, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 2)
$c.NumberFormat = '@'
$c.Value = '15:20:09'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 6)
$c.NumberFormat = '@'
$c.Value = 'Group C'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 7)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 8)
$c.NumberFormat = '@'
$c.Value = 'Шддддорпа'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 9)
$c.NumberFormat = '@'
$c.Value = 'Уруоко'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 10)
$c.NumberFormat = '@'
$c.Value = 'Уруру'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 11)
$c.NumberFormat = '@'
$c.Value = 'Уруи'
$c.ClearFormats()
$c = $ws.Cells.Item(24, 13)
$c.NumberFormat = '@'
$c.Value = 'партнер #3'
$c.ClearFormats()

# Row 25
$c = $ws.Cells.Item(25, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 2)
$c.NumberFormat = '@'
$c.Value = '15:31:37'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 6)
$c.NumberFormat = '@'
$c.Value = 'Group C'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 7)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 8)
$c.NumberFormat = '@'
$c.Value = 'Зазрззазу'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 9)
$c.NumberFormat = '@'
$c.Value = 'Дудапл'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 10)
$c.NumberFormat = '@'
$c.Value = 'Луклао'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 11)
$c.NumberFormat = '@'
$c.Value = 'Оуоапл'
$c.ClearFormats()
$c = $ws.Cells.Item(25, 13)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()

# Row 26
$c = $ws.Cells.Item(26, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 2)
$c.NumberFormat = '@'
$c.Value = '15:33:52'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 6)
$c.NumberFormat = '@'
$c.Value = 'Group B'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 7)
$c.NumberFormat = '@'
$c.Value = 'sfhdcgn'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 8)
$c.NumberFormat = '@'
$c.Value = 'sfgdgb'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 9)
$c.NumberFormat = '@'
$c.Value = 'v xcdfv'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 10)
$c.NumberFormat = '@'
$c.Value = 'dfbdxf'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 11)
$c.NumberFormat = '@'
$c.Value = 'rghfg'
$c.ClearFormats()
$c = $ws.Cells.Item(26, 13)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()

# Row 27
$c = $ws.Cells.Item(27, 1)
$c.NumberFormat = '@'
$c.Value = '2023-08-21'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 2)
$c.NumberFormat = '@'
$c.Value = '15:37:44'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 3)
$c.NumberFormat = '@'
$c.Value = 'Madina Amankeldinova'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = '@'
$c.Value = 'amankeldinovam'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 5)
$c.NumberFormat = '@'
$c.Value = 'Student'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 6)
$c.NumberFormat = '@'
$c.Value = 'Group C'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 7)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 8)
$c.NumberFormat = '@'
$c.Value = 'Цзущадал'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 9)
$c.NumberFormat = '@'
$c.Value = 'Улаллп'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 10)
$c.NumberFormat = '@'
$c.Value = 'Лулплпь'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 11)
$c.NumberFormat = '@'
$c.Value = '/start'
$c.ClearFormats()
$c = $ws.Cells.Item(27, 13)
$c.NumberFormat = '@'
$c.Value = 'Help'
$c.ClearFormats()
